$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the Spanish linking-word particles in municipality/state names ---
$ws.Range("B3").Value = "Pabellón De Arteaga"
$ws.Range("B4").Value = "Rincón De Romos"
$ws.Range("B25").Value = "Guadalupe Y Calvo"
$ws.Range("B27").Value = "Hidalgo Del Parral"
$ws.Range("A50").Value = "Ciudad De México"
$ws.Range("A72").Value = "Estado De México"
$ws.Range("B73").Value = "Atizapán De Zaragoza"
$ws.Range("B79").Value = "Ecatepec De Morelos"
$ws.Range("B80").Value = "Naucalpan De Juárez"
$ws.Range("B82").Value = "San Felipe Del Progreso"
$ws.Range("B83").Value = "San Simón De Guerrero"
$ws.Range("B86").Value = "Tlalnepantla De Baz"
$ws.Range("B95").Value = "San Francisco Del Rincón"
$ws.Range("B97").Value = "Valle De Santiago"
$ws.Range("B100").Value = "Acapulco De Juárez"
$ws.Range("B103").Value = "Chilapa De Álvarez"
$ws.Range("B104").Value = "Chilpancingo De Los Bravo"
$ws.Range("B115").Value = "Tulancingo De Bravo"
$ws.Range("B123").Value = "Cuautitlán De García Barragán"
$ws.Range("B124").Value = "Encarnación De Díaz"
$ws.Range("B128").Value = "Lagos De Moreno"
$ws.Range("B132").Value = "San Cristóbal De La Barranca"
$ws.Range("B133").Value = "San Diego De Alejandría"
$ws.Range("B134").Value = "San Juan De Los Lagos"
$ws.Range("B137").Value = "Tizapán El Alto"
$ws.Range("B138").Value = "Tlajomulco De Zúñiga"
$ws.Range("B140").Value = "Unión De Tula"
$ws.Range("B145").Value = "Zapotlán El Grande"
$ws.Range("B170").Value = "Tetela Del Volcán"
$ws.Range("B172").Value = "Amatlán De Cañas"
$ws.Range("B174").Value = "Ixtlán Del Río"
$ws.Range("B176").Value = "Santa María Del Oro"
$ws.Range("B181").Value = "Coicoyán De Las Flores"
$ws.Range("B182").Value = "Oaxaca De Juárez"
$ws.Range("B183").Value = "Putla Villa De Guerrero"
$ws.Range("B193").Value = "Los Reyes De Juárez"
$ws.Range("B196").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B201").Value = "Amealco De Bonfil"
$ws.Range("B202").Value = "San Juan Del Río"
$ws.Range("B222").Value = "Nacozari De García"
$ws.Range("B239").Value = "Muñoz De Domingo Arenas"
$ws.Range("B241").Value = "Tepetitla De Lardizábal"
$ws.Range("B250").Value = "Cosamaloapan De Carpio"
$ws.Range("B251").Value = "Martínez De La Torre"
$ws.Range("B268").Value = "Moyahua De Estrada"
$ws.Range("B269").Value = "Noria De Ángeles"
$ws.Range("B273").Value = "Teúl De González Ortega"
$ws.Range("B274").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B276").Value = "Villa De Cos"

# --- Remove trailing metadata rows (282-286); row 281 was already blank ---
$ws.Range("A281:D286").Clear()
